$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Panioty Fountain"
$ws.Range("B2").Value = "Calcutta"
$ws.Range("C2").Value = "Video"
$ws.Range("D2").Value = "22.5659° N, 88.3486° E"
$ws.Range("E2").Value = "https://www.youtube.com/watch?v=dQw4w9WgXcQ"

# Row 3
$ws.Range("A3").Value = "Hobson Jobson"
$ws.Range("B3").Value = "Tales from the past"
$ws.Range("C3").Value = "Video"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "https://www.youtube.com/watch?v=dQw4w9WgXcQ"

# Row 4
$ws.Range("A4").Value = "Duel at Alipore"
$ws.Range("B4").Value = "Calcutta"
$ws.Range("C4").Value = "Video"
$ws.Range("D4").Value = "22.5465° N, 88.3435° E"
$ws.Range("E4").Value = "https://youtu.be/dQw4w9WgXcQ"

# Row heights (rows grow to accommodate wrapped text)
$ws.Rows.Item(4).RowHeight = 60

# Column A best-fit width
$ws.Columns.Item(1).ColumnWidth = 15.14

# Final selection
$ws.Range("G4").Select()
